# Apply the lead-assessment review updates to the Test-Cases sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# First test scenario row: change Approved -> Rejected and record a reason.
$ws.Range("I2").Value = "Rejected"
$ws.Range("J2").Value = "test"

# Second test scenario row: change Approved -> Rejected and record a reason.
$ws.Range("I19").Value = "Rejected"
$ws.Range("J19").Value = "testasdsda"

# Move the active selection to reflect where the reviewer last clicked.
$ws.Range("G10").Select()
